$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# Update C2 value from 56 to 15
$ws1.Range("C2").Value = 15

# Add a new empty Sheet2 after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Update selection on Sheet1 to D10
$ws1.Activate()
$ws1.Range("D10").Select()
